# 19/12/2025: Update the list
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Connor Lutz / 1st Interview -> Eddie Powers / 2nd Interview
$ws.Range("D3").Value = "Eddie Powers"
$ws.Range("E3").Value = "2nd Interview"

# Row 4: Eddie Powers / 2nd Interview -> Jathan Prince / 3rd Interview
$ws.Range("D4").Value = "Jathan Prince"
$ws.Range("E4").Value = "3rd Interview"

# Row 5: 580/Legion Security/Sales Engineer (US)/Jathan Prince/3rd Interview
#     -> 708/Dash0/Sales Engineer (US) x 3/Sean Guillen/4th Interview
$ws.Range("A5").Value = 708
$ws.Range("B5").Value = "Dash0"
$ws.Range("C5").Value = "Sales Engineer (US) x 3"
$ws.Range("D5").Value = "Sean Guillen"
$ws.Range("E5").Value = "4th Interview"

# Row 6: 708/Dash0/Sales Engineer (US) x 3/Sean Guillen/4th Interview
#     -> 730/PointFive/PointFive SE EST/Yuval Shkedi/3rd Interview
$ws.Range("A6").Value = 730
$ws.Range("B6").Value = "PointFive"
$ws.Range("C6").Value = "PointFive SE EST"
$ws.Range("D6").Value = "Yuval Shkedi"
$ws.Range("E6").Value = "3rd Interview"

# Row 7: 730/PointFive/PointFive SE EST/Yuval Shkedi/3rd Interview
#     -> 773/CodeRabbit/Sales Engineer Bay Area/Boston/Peter Yoakum/CV Sent
$ws.Range("A7").Value = 773
$ws.Range("B7").Value = "CodeRabbit"
$ws.Range("C7").Value = "Sales Engineer Bay Area/Boston"
$ws.Range("D7").Value = "Peter Yoakum"
$ws.Range("E7").Value = "CV Sent"

# Row 8: Peter Yoakum / CV Sent -> Seth King / 3rd Interview
$ws.Range("D8").Value = "Seth King"
$ws.Range("E8").Value = "3rd Interview"

# Row 9: Seth King / 2nd Interview -> Seth Meldon / CV Sent
$ws.Range("D9").Value = "Seth Meldon"
$ws.Range("E9").Value = "CV Sent"

# Row 10 and 11 are unchanged.

# New row 12: 836/Spectro Cloud/CSE EMEA/Syed Imran/CV Sent
$ws.Range("A12").Value = 836
$ws.Range("B12").Value = "Spectro Cloud"
$ws.Range("C12").Value = "CSE EMEA"
$ws.Range("D12").Value = "Syed Imran"
$ws.Range("E12").Value = "CV Sent"

Write-Host "Applied update."
